# Edit WatchlistTestData.xlsx per commit:
# "Added new test cases in to watchlist, type ahead and profile search scenarios
#  Modified old test cases with new validations in to profiles, authoring and follow."
#
# Concretely, the visible data change in this workbook is the description text
# for the "whohas" watchlist test case (row 6, column B / DESCRIPTION) changing
# from "Get users who added  item in ther watchlist"
# to   "Get users who has  item in ther watchlist"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Watchlist")

$ws.Range("B6").Value = "Get users who has  item in ther watchlist"
